# Apply the "new .ttl from Google sheet has been generated" update:
#  - Y227 gets a date string that was previously blank
#  - A229/B229 get a new term ("vocab:1211" / "sex at birth") that was
#    previously blank, and Y229's date moves from 2023-12-15 to 2023-12-18
#  - 32 new terms are appended as rows 230-261 (dimension grows to AP261)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep all of these "date-looking" values as literal text (matches the
# source workbook, which stores them as plain inline strings rather than
# real date cells) instead of letting Excel auto-convert them to dates.
$ws.Range("Y227:Y261").NumberFormat = "@"

# --- Row 227: fill in previously empty Y227 ---
$ws.Cells.Item(227, 25).Value = "2023-12-18"

# --- Row 229: fill in previously empty A229/B229, update Y229 ---
$ws.Cells.Item(229, 1).Value = "vocab:1211"
$ws.Cells.Item(229, 2).Value = "sex at birth"
$ws.Cells.Item(229, 25).Value = "2023-12-18"

# --- New rows 230-261 ---
# Columns: A = concept id, B = preferred label, G = broader concept, Y = date.
# All other columns (C-F, H-X, Z-AP) stay empty, same as the rest of the sheet.
$newTerms = @(
    @(230, "vocab:1212", "female", "vocab:1211"),
    @(231, "vocab:1213", "male", "vocab:1211"),
    @(232, "vocab:1214", "lifestyle information variable", "vocab:1173"),
    @(233, "vocab:1215", "smoke bevaviour", "vocabl:1214"),
    @(234, "vocab:1216", "alcohol consumption", "vocabl:1214"),
    @(235, "vocab:1217", "occupation", "vocabl:1214"),
    @(236, "vocab:1218", "cosmetic use", "vocabl:1214"),
    @(237, "vocab:1219", "time activity pattern", "vocabl:1214"),
    @(238, "vocab:1220", "sociodemographic variable", "vocabl:1214"),
    @(239, "vocab:1221", "housing information", "vocabl:1214"),
    @(240, "vocab:1222", "combustion behaviour", "vocabl:1214"),
    @(241, "vocab:1223", "parity", "vocabl:1214"),
    @(242, "vocab:1224", "breastfeeding", "vocabl:1214"),
    @(243, "vocab:1225", "physical exercise", "vocabl:1214"),
    @(244, "vocab:1226", "personal hygiene", "vocabl:1214"),
    @(245, "vocab:1227", "personal behaviour", "vocabl:1214"),
    @(246, "vocab:1228", "personal information", "vocab:1173"),
    @(247, "vocab:1229", "height", "vocab:1228"),
    @(248, "vocab:1230", "weight", "vocab:1228"),
    @(249, "vocab:1231", "educational level", "vocab:1228"),
    @(250, "vocab:1232", "race/ethnicity", "vocab:1228"),
    @(251, "vocab:1233", "income", "vocab:1228"),
    @(252, "vocab:1234", "medical data/history", "vocab:1228"),
    @(253, "vocab:1235", "place/country of birth", "vocab:1228"),
    @(254, "vocab:1236", "food consumption", "vocab:1228"),
    @(255, "vocab:1237", "environmental factors information", "vocab:1173"),
    @(256, "vocab:1238", "consumption of local food/feed", "vocab:1237"),
    @(257, "vocab:1239", "urban versus non-urban", "vocab:1237"),
    @(258, "vocab:1240", "region", "vocab:1237"),
    @(259, "vocab:1241", "wheather conditions", "vocab:1237"),
    @(260, "vocab:1242", "passive smoking", "vocab:1237"),
    @(261, "vocab:1243", "industry", "vocab:1237")
)

foreach ($term in $newTerms) {
    $r = $term[0]
    $ws.Cells.Item($r, 1).Value = $term[1]
    $ws.Cells.Item($r, 2).Value = $term[2]
    $ws.Cells.Item($r, 7).Value = $term[3]
    $ws.Cells.Item($r, 25).Value = "2023-12-18"
}
